$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "'296.86"
$ws.Cells.Item(2,4).Style = 'Normal'
$ws.Cells.Item(2,5).Value = "'1.58%"
$ws.Cells.Item(2,5).Style = 'Normal'

$ws.Cells.Item(3,4).Value = "'42.04"
$ws.Cells.Item(3,4).Style = 'Normal'
$ws.Cells.Item(3,5).Value = "'3.89%"
$ws.Cells.Item(3,5).Style = 'Normal'

$ws.Cells.Item(4,4).Value = "'5.005"
$ws.Cells.Item(4,4).Style = 'Normal'
$ws.Cells.Item(4,5).Value = "'-0.16%"
$ws.Cells.Item(4,5).Style = 'Normal'

$ws.Cells.Item(5,4).Value = "'0.07530"
$ws.Cells.Item(5,4).Style = 'Normal'
$ws.Cells.Item(5,5).Value = "'2.57%"
$ws.Cells.Item(5,5).Style = 'Normal'

$ws.Cells.Item(6,4).Value = "'1.573"
$ws.Cells.Item(6,4).Style = 'Normal'
$ws.Cells.Item(6,5).Value = "'1.83%"
$ws.Cells.Item(6,5).Style = 'Normal'

$ws.Cells.Item(7,4).Value = "'0.9262"
$ws.Cells.Item(7,4).Style = 'Normal'
$ws.Cells.Item(7,5).Value = "'0.41%"
$ws.Cells.Item(7,5).Style = 'Normal'

$ws.Cells.Item(8,5).Value = "'0.05%"
$ws.Cells.Item(8,5).Style = 'Normal'

$ws.Cells.Item(9,4).Value = "'0.1194"
$ws.Cells.Item(9,4).Style = 'Normal'
$ws.Cells.Item(9,5).Value = "'-1.91%"
$ws.Cells.Item(9,5).Style = 'Normal'

$ws.Cells.Item(10,4).Value = "'0.1830"
$ws.Cells.Item(10,4).Style = 'Normal'
$ws.Cells.Item(10,5).Value = "'5.29%"
$ws.Cells.Item(10,5).Style = 'Normal'

$ws.Cells.Item(11,4).Value = "'0.08884"
$ws.Cells.Item(11,4).Style = 'Normal'
$ws.Cells.Item(11,5).Value = "'3.12%"
$ws.Cells.Item(11,5).Style = 'Normal'

$ws.Cells.Item(12,4).Value = "'0.04091"
$ws.Cells.Item(12,4).Style = 'Normal'
$ws.Cells.Item(12,5).Value = "'-4.45%"
$ws.Cells.Item(12,5).Style = 'Normal'

$ws.Cells.Item(13,4).Value = "'0.1046"
$ws.Cells.Item(13,4).Style = 'Normal'
$ws.Cells.Item(13,5).Value = "'-0.79%"
$ws.Cells.Item(13,5).Style = 'Normal'

$ws.Cells.Item(14,2).Value = 'TigerCash'
$ws.Cells.Item(14,3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Cells.Item(14,4).Value = "'0.005991"
$ws.Cells.Item(14,4).Style = 'Normal'
$ws.Cells.Item(14,5).Value = "'3.10%"
$ws.Cells.Item(14,5).Style = 'Normal'

$ws.Cells.Item(15,2).Value = 'LEO'
$ws.Cells.Item(15,3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(15,4).Value = "'3.360"
$ws.Cells.Item(15,4).Style = 'Normal'
$ws.Cells.Item(15,5).Value = "'0.60%"
$ws.Cells.Item(15,5).Style = 'Normal'

$ws.Cells.Item(16,2).Value = 'GateToken'
$ws.Cells.Item(16,3).Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Cells.Item(16,4).Value = "'4.383"
$ws.Cells.Item(16,4).Style = 'Normal'
$ws.Cells.Item(16,5).Value = "'2.03%"
$ws.Cells.Item(16,5).Style = 'Normal'

$ws.Cells.Item(17,2).Value = 'BitpandaEcosystemToken'
$ws.Cells.Item(17,3).Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Cells.Item(17,4).Value = "'0.3316"
$ws.Cells.Item(17,4).Style = 'Normal'
$ws.Cells.Item(17,5).Value = "'0.85%"
$ws.Cells.Item(17,5).Style = 'Normal'

$ws.Cells.Item(18,2).Value = 'MCDex'
$ws.Cells.Item(18,3).Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Cells.Item(18,4).Value = "'8.088"
$ws.Cells.Item(18,4).Style = 'Normal'
$ws.Cells.Item(18,5).Value = "'5.62%"
$ws.Cells.Item(18,5).Style = 'Normal'

$ws.Cells.Item(19,2).Value = 'ProBitToken'
$ws.Cells.Item(19,3).Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Cells.Item(19,4).Value = "'0.1392"
$ws.Cells.Item(19,4).Style = 'Normal'
$ws.Cells.Item(19,5).Value = "'0.05%"
$ws.Cells.Item(19,5).Style = 'Normal'

$ws.Cells.Item(20,2).Value = 'ZBToken'
$ws.Cells.Item(20,3).Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Cells.Item(20,4).Value = "'0.3301"
$ws.Cells.Item(20,4).Style = 'Normal'
$ws.Cells.Item(20,5).Value = "'20.12%"
$ws.Cells.Item(20,5).Style = 'Normal'

$ws.Cells.Item(21,2).Value = 'BitForexToken'
$ws.Cells.Item(21,3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Cells.Item(21,4).Value = "'0.001280"
$ws.Cells.Item(21,4).Style = 'Normal'
$ws.Cells.Item(21,5).Value = "'0.85%"
$ws.Cells.Item(21,5).Style = 'Normal'

$ws.Cells.Item(22,4).Value = "'0.04106"
$ws.Cells.Item(22,4).Style = 'Normal'
$ws.Cells.Item(22,5).Value = "'4.32%"
$ws.Cells.Item(22,5).Style = 'Normal'

$ws.Cells.Item(23,4).Value = "'0.001266"
$ws.Cells.Item(23,4).Style = 'Normal'
$ws.Cells.Item(23,5).Value = "'0.36%"
$ws.Cells.Item(23,5).Style = 'Normal'

$ws.Cells.Item(24,4).Value = "'0.003892"
$ws.Cells.Item(24,4).Style = 'Normal'
$ws.Cells.Item(24,5).Value = "'3.08%"
$ws.Cells.Item(24,5).Style = 'Normal'

$ws.Cells.Item(25,4).Value = "'0.0001231"
$ws.Cells.Item(25,4).Style = 'Normal'
$ws.Cells.Item(25,5).Value = "'-4.00%"
$ws.Cells.Item(25,5).Style = 'Normal'

$ws.Cells.Item(38,4).Value = "'0.02417"
$ws.Cells.Item(38,4).Style = 'Normal'
$ws.Cells.Item(38,5).Value = "'5.26%"
$ws.Cells.Item(38,5).Style = 'Normal'

$ws.Cells.Item(39,4).Value = "'0.05203"
$ws.Cells.Item(39,4).Style = 'Normal'
$ws.Cells.Item(39,5).Value = "'4.44%"
$ws.Cells.Item(39,5).Style = 'Normal'

$ws.Cells.Item(40,4).Value = "'0.006310"
$ws.Cells.Item(40,4).Style = 'Normal'
$ws.Cells.Item(40,5).Value = "'11.51%"
$ws.Cells.Item(40,5).Style = 'Normal'

$ws.Cells.Item(41,4).Value = "'0.007804"
$ws.Cells.Item(41,4).Style = 'Normal'
$ws.Cells.Item(41,5).Value = "'2.02%"
$ws.Cells.Item(41,5).Style = 'Normal'

$ws.Cells.Item(42,4).Value = "'0.1327"
$ws.Cells.Item(42,4).Style = 'Normal'
$ws.Cells.Item(42,5).Value = "'3.45%"
$ws.Cells.Item(42,5).Style = 'Normal'

$ws.Cells.Item(43,4).Value = "'0.007409"
$ws.Cells.Item(43,4).Style = 'Normal'
$ws.Cells.Item(43,5).Value = "'0.50%"
$ws.Cells.Item(43,5).Style = 'Normal'

$ws.Cells.Item(44,4).Value = "'0.007387"
$ws.Cells.Item(44,4).Style = 'Normal'
$ws.Cells.Item(44,5).Value = "'-5.09%"
$ws.Cells.Item(44,5).Style = 'Normal'

$ws.Cells.Item(45,4).Value = "'0.2951"
$ws.Cells.Item(45,4).Style = 'Normal'
$ws.Cells.Item(45,5).Value = "'-7.21%"
$ws.Cells.Item(45,5).Style = 'Normal'

$ws.Cells.Item(46,4).Value = "'0.00006436"
$ws.Cells.Item(46,4).Style = 'Normal'
$ws.Cells.Item(46,5).Value = "'1.28%"
$ws.Cells.Item(46,5).Style = 'Normal'

$ws.Cells.Item(47,4).Value = "'0.00000000751"
$ws.Cells.Item(47,4).Style = 'Normal'
$ws.Cells.Item(47,5).Value = "'0.00%"
$ws.Cells.Item(47,5).Style = 'Normal'

$ws.Cells.Item(48,4).Value = "'0.03425"
$ws.Cells.Item(48,4).Style = 'Normal'
$ws.Cells.Item(48,5).Value = "'67.43%"
$ws.Cells.Item(48,5).Style = 'Normal'

$ws.Cells.Item(49,4).Value = "'0.004207"
$ws.Cells.Item(49,4).Style = 'Normal'
$ws.Cells.Item(49,5).Value = "'0.12%"
$ws.Cells.Item(49,5).Style = 'Normal'

$ws.Cells.Item(50,4).Value = "'0.00002103"
$ws.Cells.Item(50,4).Style = 'Normal'
$ws.Cells.Item(50,5).Value = "'0.00%"
$ws.Cells.Item(50,5).Style = 'Normal'

$ws.Cells.Item(51,4).Value = "'0.0002003"
$ws.Cells.Item(51,4).Style = 'Normal'
$ws.Cells.Item(51,5).Value = "'0.00%"
$ws.Cells.Item(51,5).Style = 'Normal'
